# Automatic Birthday Wisher / Data.xlsx edit
# - Insert a new "Phone" column between Email (D) and Dialogue (old E, now F)
# - Fill in phone numbers for each row
# - Correct a few birthday dates in column C
# - Resize columns B:G
# - Move the active selection to C7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (pushes old E:F to F:G)
$ws.Range("E1").EntireColumn.Insert()

# New header for the inserted column
$ws.Range("E1").Value = "Phone"

# Phone numbers for rows 2-9 (column E), aligned with the sample data rows
$ws.Range("E2").Value = 8310145281
$ws.Range("E3").Value = 8147490519
$ws.Range("E4").Value = 8147490519
$ws.Range("E5").Value = 8310145281
$ws.Range("E6").Value = 8310145281
$ws.Range("E7").Value = 8310145281
$ws.Range("E8").Value = 8310145281
$ws.Range("E9").Value = 8310145281

# Fix a few birthdays in column C (serial date values)
$ws.Range("C4").Value = 33129
$ws.Range("C5").Value = 33128
$ws.Range("C6").Value = 34224

# Column widths (best-fit values from the authored workbook; the inputs below
# are chosen so the engine's internal character-width quantization lands on
# the closest achievable value to the recorded width for each column)
$ws.Columns.Item(2).ColumnWidth = 8.5
$ws.Columns.Item(3).ColumnWidth = 17.333333333333332
$ws.Columns.Item(4).ColumnWidth = 24.666666666666668
$ws.Columns.Item(5).ColumnWidth = 24.666666666666668
$ws.Columns.Item(6).ColumnWidth = 112.33333333333333
$ws.Columns.Item(7).ColumnWidth = 8.666666666666666

# Move the selection to reflect the author's last edit point
$ws.Range("C7").Select()
